$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1106325200512855
$ws.Range("H2").Value = 66.82167064216205
$ws.Range("I2").Value = -5.557205918961192
$ws.Range("G3").Value = 0.1045008284432751
$ws.Range("H3").Value = -11.63683356296326
$ws.Range("G4").Value = -0.3731008860193334
$ws.Range("H4").Value = -35.67534316738064
$ws.Range("G5").Value = -0.4515751981242114
$ws.Range("H5").Value = -13.17252539373433
$ws.Range("G6").Value = 0.2163282386318862
$ws.Range("H6").Value = 9.728195115900169
$ws.Range("G7").Value = 0.2924639865425127
$ws.Range("H7").Value = 41.02650768504693
$ws.Range("G8").Value = 0.1222985317607402
$ws.Range("H8").Value = 20.02287716066669
$ws.Range("G9").Value = 0.09203189303018679
$ws.Range("H9").Value = -27.2360464336911
$ws.Range("G10").Value = 0.003077960097416877
$ws.Range("H10").Value = -94.9901203159807
$ws.Range("G11").Value = 0.05477080791912636
$ws.Range("H11").Value = 9.694857150528582
$ws.Range("G12").Value = 0.109050188782545
$ws.Range("H12").Value = 17.80649681042373
$ws.Range("G13").Value = 0.08028478430678125
$ws.Range("H13").Value = 5.34993792801184
$ws.Range("G14").Value = 0.2345321142118025
$ws.Range("H14").Value = 3.792887584580119
$ws.Range("G15").Value = 0.219458702597369
$ws.Range("H15").Value = -10.92043442977264
$ws.Range("G16").Value = 0.134505459660485
$ws.Range("H16").Value = 18.25190348937222
$ws.Range("G17").Value = 0.139241079714748
$ws.Range("H17").Value = -6.805043261274515
$ws.Range("G18").Value = -0.01545694238832697
$ws.Range("H18").Value = -72.66566319772664
$ws.Range("G19").Value = -0.01459884274666243
$ws.Range("H19").Value = -160.2715291457827
$ws.Range("G20").Value = 0.1206641947765292
$ws.Range("H20").Value = 41.85580770438012
$ws.Range("G21").Value = 0.07297582058655232
$ws.Range("H21").Value = 11.49182368068946
$ws.Range("G22").Value = 0.2134265662448721
$ws.Range("H22").Value = 11.41401957428574
$ws.Range("G23").Value = 0.198639697492388
$ws.Range("H23").Value = -7.912399617681892
$ws.Range("G24").Value = -0.006690312833867494
$ws.Range("H24").Value = -75.80323672674794
$ws.Range("G25").Value = 0.004826247483038267
$ws.Range("H25").Value = 120.7515192135427
$ws.Range("G26").Value = 0.1941684885979973
$ws.Range("H26").Value = -5.221860130683438
$ws.Range("G27").Value = 0.1862750396862572
$ws.Range("H27").Value = -3.426661746769963
$ws.Range("G28").Value = 0.0600403918811617
$ws.Range("H28").Value = -10.27186359362297
$ws.Range("G29").Value = 0.09652625932949457
$ws.Range("H29").Value = 2.399291253309005